$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.716448177483187
$ws.Cells.Item(2, 3).Value = 0.1654998379516996
$ws.Cells.Item(2, 4).Value = 0.1997473665260401
$ws.Cells.Item(2, 5).Value = 0.1554444581064942
$ws.Cells.Item(2, 6).Value = 1.203585321672399
$ws.Cells.Item(2, 9).Value = 0.5308228575646865
$ws.Cells.Item(2, 10).Value = 0.1612949153980807
$ws.Cells.Item(2, 13).Value = 0.3300329971830394
$ws.Cells.Item(2, 14).Value = 1.185044160131305
$ws.Cells.Item(2, 15).Value = 2.763206647604335
$ws.Cells.Item(3, 2).Value = 0.6398855889804906
$ws.Cells.Item(3, 3).Value = 0.145388895777586
$ws.Cells.Item(3, 4).Value = 0.1980223503761849
$ws.Cells.Item(3, 5).Value = 0.15520883428999
$ws.Cells.Item(3, 6).Value = 1.202166773392378
$ws.Cells.Item(3, 9).Value = 0.5359569668228801
$ws.Cells.Item(3, 10).Value = 0.1618627829399131
$ws.Cells.Item(3, 13).Value = 0.3090634346379275
$ws.Cells.Item(3, 14).Value = 1.19461734461435
$ws.Cells.Item(3, 15).Value = 2.763768563917949
$ws.Cells.Item(4, 2).Value = 0.5928852870376033
$ws.Cells.Item(4, 3).Value = 0.1330184562545753
$ws.Cells.Item(4, 4).Value = 0.1970311934770166
$ws.Cells.Item(4, 5).Value = 0.1551299897060936
$ws.Cells.Item(4, 6).Value = 1.202008855747728
$ws.Cells.Item(4, 9).Value = 0.5394442511005373
$ws.Cells.Item(4, 10).Value = 0.162290088694089
$ws.Cells.Item(4, 13).Value = 0.2962820982477652
$ws.Cells.Item(4, 14).Value = 1.200953084107304
$ws.Cells.Item(4, 15).Value = 2.765888525533626
$ws.Cells.Item(5, 2).Value = 0.5737358378642625
$ws.Cells.Item(5, 3).Value = 0.1279720694125501
$ws.Cells.Item(5, 4).Value = 0.1966444599936565
$ws.Cells.Item(5, 5).Value = 0.1551144391353318
$ws.Cells.Item(5, 6).Value = 1.202123802510435
$ws.Cells.Item(5, 9).Value = 0.5409494834177337
$ws.Cells.Item(5, 10).Value = 0.1624839963425941
$ws.Cells.Item(5, 13).Value = 0.2910975748828619
$ws.Cells.Item(5, 14).Value = 1.203650151577342
$ws.Cells.Item(5, 15).Value = 2.767198470183587
$ws.Cells.Item(6, 2).Value = 0.5705563360137091
$ws.Cells.Item(6, 3).Value = 0.1271338063328926
$ws.Cells.Item(6, 4).Value = 0.1965812823203734
$ws.Cells.Item(6, 5).Value = 0.1551128590834594
$ws.Cells.Item(6, 6).Value = 1.20215371832758
$ws.Cells.Item(6, 9).Value = 0.5412045052602963
$ws.Cells.Item(6, 10).Value = 0.1625173891744573
$ws.Cells.Item(6, 13).Value = 0.2902381462740706
$ws.Cells.Item(6, 14).Value = 1.204104958324194
$ws.Cells.Item(6, 15).Value = 2.767442917219853
$ws.Cells.Item(7, 2).Value = 0.5926270149053039
$ws.Cells.Item(7, 3).Value = 0.1329504201617056
$ws.Cells.Item(7, 4).Value = 0.1970259082370731
$ws.Cells.Item(7, 5).Value = 0.1551297128223581
$ws.Cells.Item(7, 6).Value = 1.202009679997076
$ws.Cells.Item(7, 9).Value = 0.5394642106384282
$ws.Cells.Item(7, 10).Value = 0.1622926237174163
$ws.Cells.Item(7, 13).Value = 0.2962120804503385
$ws.Cells.Item(7, 14).Value = 1.200988991127971
$ws.Cells.Item(7, 15).Value = 2.765904386270563
$ws.Cells.Item(8, 2).Value = 0.6900481235524296
$ws.Cells.Item(8, 3).Value = 0.1585703508069685
$ws.Cells.Item(8, 4).Value = 0.1991385061598976
$ws.Cells.Item(8, 5).Value = 0.1553495672339871
$ws.Cells.Item(8, 6).Value = 1.202948197420788
$ws.Cells.Item(8, 9).Value = 0.5325235358344607
$ws.Cells.Item(8, 10).Value = 0.161474398175649
$ws.Cells.Item(8, 13).Value = 0.3227833566020166
$ws.Cells.Item(8, 14).Value = 1.18825004005302
$ws.Cells.Item(8, 15).Value = 2.763031829793562
$ws.Cells.Item(9, 2).Value = 0.881122895651913
$ws.Cells.Item(9, 3).Value = 0.2086254696966705
$ws.Cells.Item(9, 4).Value = 0.2038182534638366
$ws.Cells.Item(9, 5).Value = 0.1563020753516291
$ws.Cells.Item(9, 6).Value = 1.210448915933384
$ws.Cells.Item(9, 9).Value = 0.5215741899957109
$ws.Cells.Item(9, 10).Value = 0.1604937151413282
$ws.Cells.Item(9, 13).Value = 0.3756244571425142
$ws.Cells.Item(9, 14).Value = 1.166897193856968
$ws.Cells.Item(9, 15).Value = 2.771498834094103
$ws.Cells.Item(10, 2).Value = 1.021481722437727
$ws.Cells.Item(10, 3).Value = 0.2452792984611847
$ws.Cells.Item(10, 4).Value = 0.2075807410899841
$ws.Cells.Item(10, 5).Value = 0.157318720938239
$ws.Cells.Item(10, 6).Value = 1.219416372708011
$ws.Cells.Item(10, 9).Value = 0.5151576068202992
$ws.Cells.Item(10, 10).Value = 0.1601536291780192
$ws.Cells.Item(10, 13).Value = 0.4148835828581383
$ws.Cells.Item(10, 14).Value = 1.153415967744955
$ws.Cells.Item(10, 15).Value = 2.786344283632189
$ws.Cells.Item(11, 2).Value = 1.085321043901331
$ws.Cells.Item(11, 3).Value = 0.2619260499060658
$ws.Cells.Item(11, 4).Value = 0.2093621935763395
$ws.Cells.Item(11, 5).Value = 0.1578498232668082
$ws.Cells.Item(11, 6).Value = 1.224247958757928
$ws.Cells.Item(11, 9).Value = 0.5125932227040941
$ws.Cells.Item(11, 10).Value = 0.1600815571731289
$ws.Cells.Item(11, 13).Value = 0.432836280545061
$ws.Cells.Item(11, 14).Value = 1.147761183617668
$ws.Cells.Item(11, 15).Value = 2.794977244958176
$ws.Cells.Item(12, 2).Value = 1.109492791744287
$ws.Cells.Item(12, 3).Value = 0.2682256051931233
$ws.Cells.Item(12, 4).Value = 0.2100467691941077
$ws.Cells.Item(12, 5).Value = 0.1580607846934647
$ws.Cells.Item(12, 6).Value = 1.226185793187668
$ws.Cells.Item(12, 9).Value = 0.5116732397633683
$ws.Cells.Item(12, 10).Value = 0.1600661484132218
$ws.Cells.Item(12, 13).Value = 0.4396476769409787
$ws.Cells.Item(12, 14).Value = 1.14568851901582
$ws.Cells.Item(12, 15).Value = 2.798517045200413
$ws.Cells.Item(13, 2).Value = 1.104287119217929
$ws.Cells.Item(13, 3).Value = 0.2668690749957818
$ws.Cells.Item(13, 4).Value = 0.209898890918538
$ws.Cells.Item(13, 5).Value = 0.1580149129163608
$ws.Cells.Item(13, 6).Value = 1.225763632141707
$ws.Cells.Item(13, 9).Value = 0.5118691004029614
$ws.Cells.Item(13, 10).Value = 0.1600689384413627
$ws.Cells.Item(13, 13).Value = 0.4381801438148258
$ws.Cells.Item(13, 14).Value = 1.146131850495415
$ws.Cells.Item(13, 15).Value = 2.797742640936434
$ws.Cells.Item(14, 2).Value = 1.087309732006418
$ws.Cells.Item(14, 3).Value = 0.2624444040992557
$ws.Cells.Item(14, 4).Value = 0.209418314447305
$ws.Cells.Item(14, 5).Value = 0.1578669820276311
$ws.Cells.Item(14, 6).Value = 1.224405216460696
$ws.Cells.Item(14, 9).Value = 0.5125165105337857
$ws.Cells.Item(14, 10).Value = 0.1600800513222325
$ws.Cells.Item(14, 13).Value = 0.4333963976941817
$ws.Cells.Item(14, 14).Value = 1.147589288035945
$ws.Cells.Item(14, 15).Value = 2.795263039308168
$ws.Cells.Item(15, 2).Value = 1.076910186882287
$ws.Cells.Item(15, 3).Value = 0.2597336088314535
$ws.Cells.Item(15, 4).Value = 0.2091252448615251
$ws.Cells.Item(15, 5).Value = 0.1577776514281517
$ws.Cells.Item(15, 6).Value = 1.22358724210919
$ws.Cells.Item(15, 9).Value = 0.5129197250847142
$ws.Cells.Item(15, 10).Value = 0.1600884058616927
$ws.Cells.Item(15, 13).Value = 0.4304679106345688
$ws.Cells.Item(15, 14).Value = 1.148490953553782
$ws.Cells.Item(15, 15).Value = 2.793779474379676
$ws.Cells.Item(16, 2).Value = 1.017309348628203
$ws.Cells.Item(16, 3).Value = 0.244190820921375
$ws.Cells.Item(16, 4).Value = 0.2074657191061959
$ws.Cells.Item(16, 5).Value = 0.157285391038279
$ws.Cells.Item(16, 6).Value = 1.219115758217043
$ws.Cells.Item(16, 9).Value = 0.5153323392106337
$ws.Cells.Item(16, 10).Value = 0.160160002123547
$ws.Cells.Item(16, 13).Value = 0.4137121839296967
$ws.Cells.Item(16, 14).Value = 1.153795134015255
$ws.Cells.Item(16, 15).Value = 2.785817956708058
$ws.Cells.Item(17, 2).Value = 0.980742526955737
$ws.Cells.Item(17, 3).Value = 0.234648636950709
$ws.Cells.Item(17, 4).Value = 0.2064655012970746
$ws.Cells.Item(17, 5).Value = 0.1570009649044408
$ws.Cells.Item(17, 6).Value = 1.216565348463774
$ws.Cells.Item(17, 9).Value = 0.5169032870569481
$ws.Cells.Item(17, 10).Value = 0.1602250892583754
$ws.Cells.Item(17, 13).Value = 0.4034567935567921
$ws.Cells.Item(17, 14).Value = 1.157171446305611
$ws.Cells.Item(17, 15).Value = 2.781415534259168
$ws.Cells.Item(18, 2).Value = 0.9597093364208149
$ws.Cells.Item(18, 3).Value = 0.2291576765717025
$ws.Cells.Item(18, 4).Value = 0.2058967839044641
$ws.Cells.Item(18, 5).Value = 0.1568438312565164
$ws.Cells.Item(18, 6).Value = 1.215169222398089
$ws.Cells.Item(18, 9).Value = 0.5178402202562111
$ws.Cells.Item(18, 10).Value = 0.1602703036012514
$ws.Cells.Item(18, 13).Value = 0.3975669880894941
$ws.Cells.Item(18, 14).Value = 1.159158395587355
$ws.Cells.Item(18, 15).Value = 2.779060279741685
$ws.Cells.Item(19, 2).Value = 0.9525877405665142
$ws.Cells.Item(19, 3).Value = 0.227298102813478
$ws.Cells.Item(19, 4).Value = 0.2057053584737361
$ws.Cells.Item(19, 5).Value = 0.1567917388644275
$ws.Cells.Item(19, 6).Value = 1.214708677037876
$ws.Cells.Item(19, 9).Value = 0.5181631769639559
$ws.Cells.Item(19, 10).Value = 0.1602869482541394
$ws.Cells.Item(19, 13).Value = 0.3955743273084522
$ws.Cells.Item(19, 14).Value = 1.159838869798506
$ws.Cells.Item(19, 15).Value = 2.778293202352728
$ws.Cells.Item(20, 2).Value = 0.9846352323193628
$ws.Cells.Item(20, 3).Value = 0.2356646845793762
$ws.Cells.Item(20, 4).Value = 0.2065712955451886
$ws.Cells.Item(20, 5).Value = 0.1570305740615581
$ws.Cells.Item(20, 6).Value = 1.216829515842051
$ws.Cells.Item(20, 9).Value = 0.5167326029726667
$ws.Cells.Item(20, 10).Value = 0.1602173556298965
$ws.Cells.Item(20, 13).Value = 0.4045475866673414
$ws.Cells.Item(20, 14).Value = 1.156807376758131
$ws.Cells.Item(20, 15).Value = 2.781865868085788
$ws.Cells.Item(21, 2).Value = 1.092296488922841
$ws.Cells.Item(21, 3).Value = 0.2637441537028451
$ws.Cells.Item(21, 4).Value = 0.2095592011590099
$ws.Cells.Item(21, 5).Value = 0.1579101659275466
$ws.Cells.Item(21, 6).Value = 1.22480127871809
$ws.Cells.Item(21, 9).Value = 0.5123249629452182
$ws.Cells.Item(21, 10).Value = 0.160076464684721
$ws.Cells.Item(21, 13).Value = 0.4348011471578488
$ws.Cells.Item(21, 14).Value = 1.147159339667773
$ws.Cells.Item(21, 15).Value = 2.795984009016337
$ws.Cells.Item(22, 2).Value = 1.162642322198735
$ws.Cells.Item(22, 3).Value = 0.2820710478123374
$ws.Cells.Item(22, 4).Value = 0.2115700984430191
$ws.Cells.Item(22, 5).Value = 0.1585424003534897
$ws.Cells.Item(22, 6).Value = 1.230642082548584
$ws.Cells.Item(22, 9).Value = 0.5097421791049577
$ws.Cells.Item(22, 10).Value = 0.1600536498057821
$ws.Cells.Item(22, 13).Value = 0.4546497892959707
$ws.Cells.Item(22, 14).Value = 1.141254089850278
$ws.Cells.Item(22, 15).Value = 2.806788968843165
$ws.Cells.Item(23, 2).Value = 1.125099410752227
$ws.Cells.Item(23, 3).Value = 0.2722919912870339
$ws.Cells.Item(23, 4).Value = 0.2104915491313903
$ws.Cells.Item(23, 5).Value = 0.1581997233864136
$ws.Cells.Item(23, 6).Value = 1.227467002841905
$ws.Cells.Item(23, 9).Value = 0.5110933687503731
$ws.Cells.Item(23, 10).Value = 0.1600594885475246
$ws.Cells.Item(23, 13).Value = 0.4440493392536595
$ws.Cells.Item(23, 14).Value = 1.144369217799529
$ws.Cells.Item(23, 15).Value = 2.800877649405749
$ws.Cells.Item(24, 2).Value = 0.9828753726382615
$ws.Cells.Item(24, 3).Value = 0.2352053450990184
$ws.Cells.Item(24, 4).Value = 0.206523446271234
$ws.Cells.Item(24, 5).Value = 0.1570171678661652
$ws.Cells.Item(24, 6).Value = 1.216709867295975
$ws.Cells.Item(24, 9).Value = 0.5168096640494326
$ws.Cells.Item(24, 10).Value = 0.1602208277230659
$ws.Cells.Item(24, 13).Value = 0.404054419903737
$ws.Cells.Item(24, 14).Value = 1.156971829736904
$ws.Cells.Item(24, 15).Value = 2.781661724716372
$ws.Cells.Item(25, 2).Value = 0.8294331560998671
$ws.Cells.Item(25, 3).Value = 0.1951049809730421
$ws.Cells.Item(25, 4).Value = 0.2024950905911851
$ws.Cells.Item(25, 5).Value = 0.1559886840671041
$ws.Cells.Item(25, 6).Value = 1.207813176399512
$ws.Cells.Item(25, 9).Value = 0.524250798995098
$ws.Cells.Item(25, 10).Value = 0.1606922174406762
$ws.Cells.Item(25, 13).Value = 0.361252026773208
$ws.Cells.Item(25, 14).Value = 1.172285828940957
$ws.Cells.Item(25, 15).Value = 2.767695776451205
